$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source workbook lists football match results/odds. Two pairs of rows
# (87/88 and 111/112) had their match records swapped (the row-sequence
# column A keeps its original 85/86/109/110 numbering, but every other
# field - id, teams, odds, etc. - belongs to the other row of the pair).
# Fix this by swapping the full B:AC payload between each pair of rows.

function Swap-Rows([int]$rowA, [int]$rowB) {
    $rangeA = $ws.Range("B$rowA`:AC$rowA")
    $rangeB = $ws.Range("B$rowB`:AC$rowB")

    $valuesA = $rangeA.Value()
    $valuesB = $rangeB.Value()

    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}

Swap-Rows 87 88
Swap-Rows 111 112
